$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5: "Data Analytics" (was previously row 6), B stays 200, C becomes 232
$ws.Range("A5").Value = "Data Analytics"
$ws.Range("B5").Value = 200
$ws.Range("C5").Value = 232

# Row 6: "Data Mining" (was previously row 7), B/C unchanged (200 / 0)
$ws.Range("A6").Value = "Data Mining"
$ws.Range("B6").Value = 200
$ws.Range("C6").Value = 0

# Row 7: "Data Warehousing" (was previously row 8), B/C unchanged (200 / 0)
$ws.Range("A7").Value = "Data Warehousing"
$ws.Range("B7").Value = 200
$ws.Range("C7").Value = 0

# Row 8: "The Importance of Big Data" (was previously row 5), B becomes 300, C stays 0
$ws.Range("A8").Value = "The Importance of Big Data"
$ws.Range("B8").Value = 300
$ws.Range("C8").Value = 0
